{"js": "// Insert a new paragraph \"git reset (make all files untracked)\" right\n// after the \"git push \u2013u origin master\" paragraph, with the document's\n// lone \"_GoBack\" bookmark moved so it now sits inside the new paragraph\n// (between \"...untracked\" and the closing \")\").\n\n// 1) Remove the existing \"_GoBack\" bookmark first (it currently lives in\n//    the last, otherwise-empty paragraph of the document). Doing this\n//    before inserting the new one avoids any ambiguity about which\n//    same-named bookmark a later lookup would resolve to.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// 2) Find the paragraph that holds exactly \"git push \u2013u origin master\"\n//    (the text also appears, as a substring, inside an earlier\n//    paragraph, so match the full paragraph text rather than searching).\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst targetText = \"git push \\u2013u origin master\";\nlet anchorParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  if (paragraphs.items[i].text === targetText) {\n    anchorParagraph = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchorParagraph) {\n  throw new Error('Could not find paragraph \"' + targetText + '\"');\n}\n\n// 3) Insert the new paragraph right after it; inserting via the\n//    paragraph inherits its paragraph formatting (spacing/theme fonts),\n//    matching the surrounding paragraphs.\nconst newParagraph = anchorParagraph.insertParagraph(\n  \"git reset (make all files untracked)\",\n  \"After\"\n);\nawait context.sync();\n\n// 4) Split the new paragraph's text so the bookmark sits right after\n//    \"...untracked\" and before the closing \")\".\nconst hits = newParagraph.search(\"untracked\");\nhits.load(\"text\");\nawait context.sync();\n\nif (hits.items.length === 0) {\n  throw new Error(\"Could not find inserted text to place the bookmark\");\n}\n\nconst bookmarkPosition = hits.items[0].getRange(\"End\");\nbookmarkPosition.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# Insert a new paragraph \"git reset (make all files untracked)\" right\n# after the \"git push \u2013u origin master\" paragraph, with the document's\n# lone \"_GoBack\" bookmark moved so it now sits inside the new paragraph\n# (between \"...untracked\" and the closing \")\").\n\n$d = $word.ActiveDocument\n\n# 1) Locate the paragraph whose text is exactly \"git push \u2013u origin\n#    master\" (the same words also appear, as a substring, inside an\n#    earlier paragraph, so match the whole paragraph text rather than\n#    just searching for the phrase).\n$targetText = \"git push \" + [char]0x2013 + \"u origin master\"\n$anchorIndex = 0\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)\n    if ($t -eq $targetText) {\n        $anchorIndex = $i\n    }\n}\n\nif ($anchorIndex -eq 0) {\n    Write-Output \"Could not find anchor paragraph\"\n} else {\n    # 2) Insert a new paragraph right after it; this inherits the\n    #    anchor paragraph's formatting (spacing/theme fonts), matching\n    #    the surrounding paragraphs.\n    $d.Paragraphs.Item($anchorIndex).Range.InsertParagraphAfter()\n    $newParagraph = $d.Paragraphs.Item($anchorIndex + 1)\n    $newParagraph.Range.Text = \"git reset (make all files untracked)\"\n\n    # 3) Split the new paragraph's text so the bookmark sits right\n    #    after \"...untracked\" and before the closing \")\". Adding a\n    #    bookmark with an already-used name moves it (Word bookmark\n    #    names are unique), so this both repositions and removes the\n    #    original \"_GoBack\" bookmark in one step.\n    $newRange = $newParagraph.Range\n    $newText = $newRange.Text\n    $offset = $newRange.Start + $newText.IndexOf(\"untracked\") + \"untracked\".Length\n    $bookmarkRange = $d.Range($offset, $offset)\n    $d.Bookmarks.Add(\"_GoBack\", $bookmarkRange)\n}\n"}
